$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ubicacion")

$ws.Range("A4").Value = "unicentro"
$ws.Range("B4").Value = "unnicentro"
$ws.Range("C4").Value = "Centro Comercial"
$ws.Range("D4").Value = "Diurno"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "hala"
$ws.Range("G4").Value = "'12124"
$ws.Range("H4").Value = "soymejorqueventura"
